$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 0.4571641683578491
$ws.Cells.Item(1, 2).Value = 0.9303494095802307
$ws.Cells.Item(1, 3).Value = 1.03879976272583
$ws.Cells.Item(1, 4).Value = 5.127120494842529
$ws.Cells.Item(1, 5).Value = 1.270225882530212
